# Conserto do erro com o rotulo da coluna 2050 nas tabelas e retirada das
# linhas com total das tabelas.
#
# For each of the first five tables (worksheets 1-5) the header cell in
# column E (row 1) incorrectly held a stray numeric value; it should hold
# the text label for the "2050" column (or "2041-2050" on the incremental
# power sheet) just like the other header cells in that row. The sixth
# worksheet has no such header cell. In addition, every worksheet's final
# "Total" row is removed.

$wb = $excel.ActiveWorkbook

# Writes $text into ($row, $col) as a genuine text value (matching the
# style/type of neighboring text header cells) without relying on
# Range.Value's "looks like a number -> becomes a Double" coercion, and
# without leaving a stray quote-prefixed style behind. We do this by
# stashing a string *formula* result (which already carries the "text"
# cell type) in a scratch cell far outside the used range, copying just
# that cell's value onto the destination, then wiping the scratch cell.
function Set-TextLabel {
    param($ws, [int]$row, [int]$col, [string]$text)

    $scratchRow = 1000
    $scratchCol = 1000
    $scratch = $ws.Cells.Item($scratchRow, $scratchCol)

    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues

    $ws.Parent.Application.CutCopyMode = 0
    $scratch.ClearContents()
}

# --- Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio
#     (MWMed)", "Atendimento a Ponta(MW)" — header E1 "2040"->"2050" label,
#     drop the Total row (row 13). ---------------------------------------
foreach ($name in @(
        "Potencia Acumulada - SIN (MW)",
        "Geracao Periodo Medio (MWMed)",
        "Atendimento a Ponta(MW)"
    )) {
    $ws = $wb.Worksheets.Item($name)
    Set-TextLabel $ws 1 5 "2050"
    $ws.Rows.Item(13).Delete()
}

# --- Sheet 4: "Potencia Incremental - SIN(MW)" — header E1 uses the
#     "2031-2040"->"2041-2050" range label, drop the Total row (row 13). ---
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
Set-TextLabel $ws4 1 5 "2041-2050"
$ws4.Rows.Item(13).Delete()

# --- Sheet 5: "Emissoes Totais (MtCO2eq)" — header E1 "2040"->"2050"
#     label; this table has no Total row to remove. ------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
Set-TextLabel $ws5 1 5 "2050"

# --- Sheet 6: "Custo Total (bilhões de R$)" — no header-label column to
#     fix, just drop the Total row (row 4). --------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
